$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns whose values are swapped between row 2 and row 3: A, Q, R, AC are
# "naturally typed" (numeric stays numeric, non-numeric text stays text).
# Column I holds numeric-looking text ("15" / "100") that must stay text.
$plainCols = @("A", "Q", "R", "AC")

foreach ($col in $plainCols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $val2 = $ws.Range($addr2).Value2
    $val3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $val3
    $ws.Range($addr3).Value2 = $val2
}

# Column I ("Antal") stores digit-only text ("15"/"100"); Value2 auto-coerces
# digit strings to numbers, so force Text format while assigning, then
# restore the default style so no stray formatting is left behind.
$i2 = $ws.Range("I2")
$i3 = $ws.Range("I3")
$iVal2 = $i2.Value2
$iVal3 = $i3.Value2

$i2.NumberFormat = "@"
$i2.Value2 = "$iVal3"
$i2.Style = "Normal"

$i3.NumberFormat = "@"
$i3.Value2 = "$iVal2"
$i3.Style = "Normal"
